# Updated cryptos list with GitHub Actions
# Applies the price/volume(1h) refresh plus a few rank re-shuffles
# described by the authoritative diff.
#
# Notes on this runtime's quirks (discovered empirically):
#  1. Named parameters (e.g. "-Row 2 -D foo") are NOT supported by this
#     PowerShell-style interpreter; Set-Row below takes its arguments
#     purely positionally: Row, B, C, D, E. Pass $null for any column
#     that should be left untouched.
#  2. Assigning a numeric-looking string straight to Range/Cells.Value
#     makes Excel silently coerce it to a real number (losing
#     trailing zeros like "89.50" -> 89.5, or turning tiny decimals
#     like "0.0000277" into scientific notation). To preserve the
#     exact original text (matching the workbook's inlineStr cells),
#     we force the cell to Text format ("@") before assigning the
#     value, then reset the cell style back to Normal afterwards so
#     we don't leave a stray explicit number-format style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Cell, $V)
    $Cell.NumberFormat = "@"
    $Cell.Value = $V
    $Cell.Style = "Normal"
}

function Set-Row {
    param($Row, $B, $C, $D, $E)
    if ($null -ne $B) { Set-CellText $ws.Cells.Item($Row, 2) $B }
    if ($null -ne $C) { Set-CellText $ws.Cells.Item($Row, 3) $C }
    if ($null -ne $D) { Set-CellText $ws.Cells.Item($Row, 4) $D }
    if ($null -ne $E) { Set-CellText $ws.Cells.Item($Row, 5) $E }
}

# Row 2: Bitcoin
Set-Row 2 $null $null "65.569.76" "  +1.06%  "

# Row 3: Ethereum
Set-Row 3 $null $null "3.397.16" "  +0.15%  "

# Row 4: TetherUSD
Set-Row 4 $null $null $null "  -0.06%  "

# Row 5: BNB
Set-Row 5 $null $null "560.47" "  -0.01%  "

# Row 6: Solana
Set-Row 6 $null $null "175.98" "  +0.64%  "

# Row 7: XRP
Set-Row 7 $null $null $null "  +0.78%  "

# Row 8: LidoStakedEther
Set-Row 8 $null $null "3.390.52" "  +0.30%  "

# Row 9: USDC
Set-Row 9 $null $null $null "  -0.08%  "

# Row 10: Dogecoin
Set-Row 10 $null $null "0.174" "  +4.93%  "

# Row 11: Cardano
Set-Row 11 $null $null "0.636" "  +0.52%  "

# Row 12: Avalanche
Set-Row 12 $null $null "53.52" "  -1.96%  "

# Row 13: ShibaInu
Set-Row 13 $null $null "0.0000277" "  +0.22%  "

# Row 14: Polkadot
Set-Row 14 $null $null "9.21" "  +0.75%  "

# Row 15: WrappedliquidstakedEther2.0
Set-Row 15 $null $null "3.939.38" "  +0.19%  "

# Row 16: Chainlink
Set-Row 16 $null $null "18.28" "  +0.01%  "

# Row 17: was WrappedEther -> now TRON
Set-Row 17 "TRON" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx" "0.119" "  +0.91%  "

# Row 18: was TRON -> now WrappedEther
Set-Row 18 "WrappedEther" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" "3.378.43" "  -0.31%  "

# Row 19: WrappedBTC
Set-Row 19 $null $null "65.467.28" "  +1.00%  "

# Row 20: Uniswap
Set-Row 20 $null $null "11.84" "  -0.54%  "

# Row 21: Polygon
Set-Row 21 $null $null $null "  +0.34%  "

# Row 22: BitcoinCash
Set-Row 22 $null $null "478.19" "  +1.78%  "

# Row 23: Toncoin
Set-Row 23 $null $null "4.93" "  -1.35%  "

# Row 24: was InternetComputer(DFINITY) -> now Litecoin
Set-Row 24 "Litecoin" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc" "89.50" "  +3.63%  "

# Row 25: was PancakeSwap -> now InternetComputer(DFINITY)
Set-Row 25 "InternetComputer(DFINITY)" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp" "14.29" "  +4.59%  "

# Row 26: was Litecoin -> now PancakeSwap
Set-Row 26 "PancakeSwap" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake" "4.10" "  -0.74%  "

# Row 27: ImmutableX
Set-Row 27 $null $null $null "  +2.09%  "

# Row 28: RenderToken
Set-Row 28 $null $null $null "  -1.79%  "

# Row 29: Filecoin
Set-Row 29 $null $null "8.73" "  -1.33%  "

# Row 30: EthereumClassic
Set-Row 30 $null $null "31.24" "  +1.91%  "

# Row 31: NEARProtocol
Set-Row 31 $null $null "6.56" "  -2.27%  "

# Row 32: Cosmos
Set-Row 32 $null $null "11.50" "  -0.50%  "

# Row 33: OKB
Set-Row 33 $null $null "62.96" "  +4.77%  "

# Row 34: Bittensor
Set-Row 34 $null $null "575.38" "  -0.71%  "

# Row 35: Hedera
Set-Row 35 $null $null $null "  -0.99%  "

# Row 36: Dai
Set-Row 36 $null $null $null "  -0.03%  "

# Row 37: Stacks
Set-Row 37 $null $null "3.68" "  +5.71%  "

# Row 38: Kaspa
Set-Row 38 $null $null $null "  +0.00%  "

# Row 39: InjectiveProtocol
Set-Row 39 $null $null "35.83" "  -0.35%  "

# Row 40: TheGraph
Set-Row 40 $null $null "0.374" "  +0.34%  "

# Row 41: PEPE
Set-Row 41 $null $null $null "  -1.91%  "

# Row 42: Maker
Set-Row 42 $null $null "3.089.77" "  -0.38%  "

# Row 43: ThetaToken
Set-Row 43 $null $null "2.80" "  -2.71%  "

# Row 44: VeChain
Set-Row 44 $null $null "0.0417" "  +0.67%  "

# Row 45: Stellar
Set-Row 45 $null $null "0.134" "  +0.32%  "

# Row 46: ApeXProtocol
Set-Row 46 $null $null $null "  -0.96%  "

# Row 47: Fetch.AI
Set-Row 47 $null $null $null "  -3.70%  "

# Row 48: FirstDigitalUSD
Set-Row 48 $null $null "1.00" "  +0.01%  "

# Row 49: Monero
Set-Row 49 $null $null "140.19" "  +1.82%  "

# Row 50: WEMIXToken
Set-Row 50 $null $null "2.56" "  -0.16%  "

# Row 51: THORChain
Set-Row 51 $null $null "8.43" "  +0.63%  "
